$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking price strings
# (e.g. "0.9995", "27.187.51") are preserved exactly as authored.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("D2").Value = "27.187.51"
$ws.Range("E2").Value = "  +0.09%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.831.58"
$ws.Range("E3").Value = "  +0.80%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.86%  "

# Row 5: BNB
$ws.Range("D5").Value = "312.62"
$ws.Range("E5").Value = "  -0.47%  "

# Row 6: USDC
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.76%  "

# Row 7: XRP
$ws.Range("D7").Value = "0.4537"
$ws.Range("E7").Value = "  +6.83%  "

# Row 8: Cardano
$ws.Range("D8").Value = "0.3755"
$ws.Range("E8").Value = "  +2.62%  "

# Row 9: Dogecoin
$ws.Range("D9").Value = "0.07341"
$ws.Range("E9").Value = "  +2.20%  "

# Row 10: Polygon
$ws.Range("D10").Value = "0.8595"
$ws.Range("E10").Value = "  +0.20%  "

# Row 11: Solana
$ws.Range("D11").Value = "21.04"
$ws.Range("E11").Value = "  +0.50%  "

# Row 12: WrappedEther
$ws.Range("D12").Value = "1.836.71"
$ws.Range("E12").Value = "  +1.01%  "

# Row 13: Chainlink
$ws.Range("D13").Value = "6.715"
$ws.Range("E13").Value = "  +1.40%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "5.354"
$ws.Range("E14").Value = "  +1.05%  "

# Row 15: Litecoin
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "92.49"
$ws.Range("E15").Value = "  +5.53%  "

# Row 16: TRON
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "0.07110"
$ws.Range("E16").Value = "  +0.73%  "

# Row 17: BinanceUSD
$ws.Range("E17").Value = "  -0.85%  "

# Row 18: ShibaInu
$ws.Range("D18").Value = "0.000008848"
$ws.Range("E18").Value = "  +0.31%  "

# Row 19: Dai
$ws.Range("D19").Value = "0.9989"
$ws.Range("E19").Value = "  -0.90%  "

# Row 20: Avalanche
$ws.Range("D20").Value = "15.04"
$ws.Range("E20").Value = "  +0.11%  "

# Row 21: WrappedBTC
$ws.Range("D21").Value = "27.173.85"
$ws.Range("E21").Value = "  -0.10%  "

# Row 22: Uniswap
$ws.Range("D22").Value = "5.203"
$ws.Range("E22").Value = "  +2.07%  "

# Row 23: Cosmos
$ws.Range("D23").Value = "11.01"
$ws.Range("E23").Value = "  +1.67%  "

# Row 24: Toncoin
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "2.002"
$ws.Range("E24").Value = "  -0.09%  "

# Row 25: Monero
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "152.20"
$ws.Range("E25").Value = "  -0.33%  "

# Row 26: LidoDAOToken
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "2.253"
$ws.Range("E26").Value = "  +6.37%  "

# Row 27: EthereumClassic
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "18.54"
$ws.Range("E27").Value = "  +1.27%  "

# Row 28: InternetComputer(DFINITY)
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "5.294"
$ws.Range("E28").Value = "  +1.24%  "

# Row 29: BitcoinCash
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "117.47"
$ws.Range("E29").Value = "  +1.91%  "

# Row 30: Stellar
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "0.08852"
$ws.Range("E30").Value = "  -0.01%  "

# Row 31: ARBITRUM
$ws.Range("B31").Value = "ARBITRUM"
$ws.Range("C31").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D31").Value = "1.200"
$ws.Range("E31").Value = "  +0.61%  "

# Row 32: ImmutableX
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "0.7513"
$ws.Range("E32").Value = "  -0.20%  "

# Row 33: HuobiToken
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "2.977"
$ws.Range("E33").Value = "  +5.12%  "

# Row 34: Filecoin
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "4.478"
$ws.Range("E34").Value = "  +0.45%  "

# Row 35: Frax
$ws.Range("B35").Value = "Frax"
$ws.Range("C35").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D35").Value = "0.9987"
$ws.Range("E35").Value = "  -0.86%  "

# Row 36: TrustWalletToken
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "1.101"
$ws.Range("E36").Value = "  -0.76%  "

# Row 37: VeChain
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.01971"
$ws.Range("E37").Value = "  +1.05%  "

# Row 38: Hedera
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.05290"
$ws.Range("E38").Value = "  +1.39%  "

# Row 39: FraxShare
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "7.291"
$ws.Range("E39").Value = "  +3.12%  "

# Row 40: TheSandbox
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.5370"
$ws.Range("E40").Value = "  +7.29%  "

# Row 41: MXToken
$ws.Range("D41").Value = "2.894"
$ws.Range("E41").Value = "  +0.56%  "

# Row 42: Algorand
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.1712"
$ws.Range("E42").Value = "  +2.64%  "

# Row 43: Aptos
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "8.658"
$ws.Range("E43").Value = "  +1.38%  "

# Row 44: Decentraland
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "0.5101"
$ws.Range("E44").Value = "  +8.80%  "

# Row 45: EnergySwap
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "10.69"
$ws.Range("E45").Value = "  +1.99%  "

# Row 46: RenderToken
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "1.972"
$ws.Range("E46").Value = "  +9.61%  "

# Row 47: Quant
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "106.25"
$ws.Range("E47").Value = "  +0.16%  "

# Row 48: NEARProtocol
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.676"
$ws.Range("E48").Value = "  +1.34%  "

# Row 49: PaxDollar
$ws.Range("D49").Value = "0.9986"
$ws.Range("E49").Value = "  -0.84%  "

# Row 50: Cronos
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.06353"
$ws.Range("E50").Value = "  -0.57%  "

# Row 51: ThetaToken
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "0.9232"
$ws.Range("E51").Value = "  +1.81%  "

